$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list: refresh Price (D) / Volume(1h) (E) for changed rows,
# and bump Hora (G) from 4 to 5 for every data row (2-51).
$updates = @{
    2 = @{ D="311.09"; E="-0.68%" }
    3 = @{ D="37.68"; E="-0.59%" }
    4 = @{ D="5.165"; E="1.90%" }
    5 = @{ D="0.07932"; E="2.17%" }
    6 = @{ D="1.922"; E="1.20%" }
    7 = @{ D="8.286"; E="1.26%" }
    8 = @{ D="2.989"; E="-3.56%" }
    9 = @{ D="0.9309"; E="1.43%" }
    10 = @{ D="0.1084"; E="-12.39%" }
    11 = @{ D="0.1925"; E="1.50%" }
    12 = @{ D="0.09157"; E="3.58%" }
    13 = @{ D="0.03298"; E="-2.73%" }
    14 = @{ D="0.09599"; E="-1.03%" }
    15 = @{ D="0.001378"; E="0.94%" }
    16 = @{ D="0.005815"; E="-1.92%" }
    17 = @{ D="3.596"; E="1.73%" }
    18 = @{ D="4.436" }
    19 = @{ D="0.3409"; E="0.00%" }
    20 = @{ D="6.419"; E="27.70%" }
    21 = @{ E="-0.65%" }
    22 = @{ E="-0.01%" }
    23 = @{ D="0.04410"; E="0.14%" }
    24 = @{ D="0.001233"; E="1.70%" }
    25 = @{ D="0.004627"; E="8.92%" }
    26 = @{  }
    27 = @{ D="0.0003991" }
    28 = @{  }
    29 = @{  }
    30 = @{  }
    31 = @{  }
    32 = @{  }
    33 = @{  }
    34 = @{  }
    35 = @{  }
    36 = @{  }
    37 = @{  }
    38 = @{  }
    39 = @{ D="0.02245"; E="4.65%" }
    40 = @{ D="0.05090"; E="2.20%" }
    41 = @{ D="0.007477"; E="-4.26%" }
    42 = @{ D="0.008936"; E="-10.26%" }
    43 = @{ D="0.1355"; E="0.68%" }
    44 = @{  }
    45 = @{ D="0.008611"; E="-11.06%" }
    46 = @{ D="0.00006624"; E="1.66%" }
    47 = @{ D="0.00000000750"; E="-0.01%" }
    48 = @{ E="-10.56%" }
    49 = @{ D="0.001000"; E="-40.76%" }
    50 = @{ D="0.00002101"; E="-0.01%" }
    51 = @{ D="0.0002001"; E="-0.01%" }
}

foreach ($rowKey in $updates.Keys) {
    $row = [int]$rowKey
    $vals = $updates[$rowKey]
    if ($vals.ContainsKey("D")) {
        $ws.Cells.Item($row, 4).Value = "'" + $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = "'" + $vals["E"]
    }
    $ws.Cells.Item($row, 7).Value = "'5"
}

Write-Output "Updated rows 2-51: prices/volumes refreshed, Hora set to 5."